$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A25").Value = "Plexiglass plate, 5mm thick, with cutouts"
$ws.Range("B25").Value = "this repository, ideally commission at your instituions workshop"
$ws.Range("D25").Value = "1"

$ws.Range("E16").Select()
